# Update to project status
# ---------------------------------------------------------------------------
# This script reproduces, via Excel COM automation, the edits described by
# the commit "Update to project status":
#   1. Mark two existing Estimates rows (Story rows 139 & 152) as Completed.
#   2. Add two brand-new Story rows to Table1 on the Estimates sheet:
#        - "Tech Debt" / "Make EmailQueue into a Background Service"
#        - "New Stories" / "Apply String Localization"
#      both marked Completed, which pushes the summary/stat block down.
#   3. Append a new data point (week 65 / 524 remaining hours) to the
#      "Burn Down" sheet and extend the burn-down chart series ranges.
#   4. Switch the active/selected sheet from "Estimates" to "Burn Down" and
#      update each sheet's remembered selection.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Flip "Completed" on two already-existing stories (rows 139 and 152).
#    Their Completed Points / Completed Hours formulas already exist on the
#    sheet, so simply setting Completed = TRUE lets them recalc on their own.
# ---------------------------------------------------------------------------
$wsEstimates = $wb.Worksheets.Item("Estimates")

$wsEstimates.Range("E139").Value = $true
$wsEstimates.Range("E152").Value = $true

# ---------------------------------------------------------------------------
# 2. Insert two new rows right after the last table row (row 180) so that
#    everything below (the gap rows + summary/statistics block) shifts down
#    by two rows automatically, including every formula that references it.
# ---------------------------------------------------------------------------
$wsEstimates.Rows.Item(181).Resize(2).Insert()

# Row 181: Tech Debt / Make EmailQueue into a Background Service
$wsEstimates.Range("A181").Value = "Tech Debt"
$wsEstimates.Range("B181").Value = "Make EmailQueue into a Background Service"
$wsEstimates.Range("C181").Value = 1
$wsEstimates.Range("D181").Formula = "=VLOOKUP(C181,Points!`$A`$1:`$C`$6,3,FALSE)"
$wsEstimates.Range("E181").Value = $true
$wsEstimates.Range("F181").Formula = "=IF(Table1[[#This Row],[Completed]],Table1[[#This Row],[Points]],0)"
$wsEstimates.Range("G181").Formula = "=IF(Table1[[#This Row],[Completed]],Table1[[#This Row],[Estimated Hours]],0)"

# Row 182: New Stories / Apply String Localization
$wsEstimates.Range("A182").Value = "New Stories"
$wsEstimates.Range("B182").Value = "Apply String Localization"
$wsEstimates.Range("C182").Value = 3
$wsEstimates.Range("D182").Formula = "=VLOOKUP(C182,Points!`$A`$1:`$C`$6,3,FALSE)"
$wsEstimates.Range("E182").Value = $true
$wsEstimates.Range("F182").Formula = "=IF(Table1[[#This Row],[Completed]],Table1[[#This Row],[Points]],0)"
$wsEstimates.Range("G182").Formula = "=IF(Table1[[#This Row],[Completed]],Table1[[#This Row],[Estimated Hours]],0)"

# Grow Table1 (and its auto filter) so the two new rows are part of it.
$loTable1 = $wsEstimates.ListObjects.Item("Table1")
$loTable1.Resize($wsEstimates.Range("A1:G182"))

# ---------------------------------------------------------------------------
# 3. Burn Down sheet: append week 65 / 524 remaining hours as row 10.
# ---------------------------------------------------------------------------
$wsBurnDown = $wb.Worksheets.Item("Burn Down")

$wsBurnDown.Range("A10").Formula = "=A9+7"
$wsBurnDown.Range("B10").Formula = "=B9+1"
$wsBurnDown.Range("C10").Value = 524

# Extend the Burn Down chart's category/value series to include the new row.
$coChart = $wsBurnDown.ChartObjects().Item(1)
$chart = $coChart.Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES('Burn Down'!`$C`$1,'Burn Down'!`$B`$2:`$B`$10,'Burn Down'!`$C`$2:`$C`$10,1)"

# ---------------------------------------------------------------------------
# 4. Update the active sheet / selections to match the saved workbook state:
#    "Burn Down" becomes the active tab, with its own remembered selection,
#    while "Estimates" keeps a (non-active) selection of its own.
# ---------------------------------------------------------------------------
$wsEstimates.Activate()
$wsEstimates.Range("E191").Select()

$wsBurnDown.Activate()
$wsBurnDown.Range("C23").Select()
